$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: Qualcomm Atheros AR9580 Wireless Network Adapter - 10.1.10.5
$ws.Range("A3").Value = "Qualcomm Atheros AR9580 Wireless Network Adapter - 10.1.10.5"
$ws.Range("B3").Value = 3
$ws.Range("C3").Value = 845
$ws.Range("D3").Value = 95

# Row 4: Qualcomm Atheros AR9580 Wireless Network Adapter - 3.0.2.201 (name unchanged)
$ws.Range("B4").Value = 3
$ws.Range("C4").Value = 1083
$ws.Range("D4").Value = 95.09999999999999

# Row 5: Realtek RTL8852AE WiFi 6 802.11ax PCIe Adapter - 6001.10.356.0
$ws.Range("A5").Value = "Realtek RTL8852AE WiFi 6 802.11ax PCIe Adapter - 6001.10.356.0"
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = 14
$ws.Range("D5").Value = 98.7

# Row 6: Totals
$ws.Range("B6").Value = 7
$ws.Range("C6").Value = 1942
